$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

for ($r = 2; $r -le 11; $r++) {
    $ws.Range("B$r").Formula = "=SUMIFS(lunch_status_quantile!`$C:`$C, lunch_status_quantile!`$A:`$A, Sheet1!`$A$r, lunch_status_quantile!`$B:`$B, 2023)*100"
    $ws.Range("C$r").Formula = "=SUMIFS(lunch_status_quantile!`$C:`$C, lunch_status_quantile!`$A:`$A, Sheet1!`$A$r, lunch_status_quantile!`$B:`$B, 2024)*100"
}

$ws.Range("E31:G34").Select()

$wb.Save()
